$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some updated Price values are numeric-looking ("310.88", "1.007", ...).
# The source column stores Price as text, so force Text format on those
# specific cells first to stop Excel from auto-converting them to numbers.
$textCells = @("D5","D7","D8","D9","D10","D11","D12","D14","D15","D16","D17","D18","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D35","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.887.58'
$ws.Range("E2").Value = '  -1.74%  '
$ws.Range("D3").Value = '1.826.20'
$ws.Range("E3").Value = '  -1.66%  '
$ws.Range("E4").Value = '  +0.56%  '
$ws.Range("D5").Value = '310.88'
$ws.Range("E5").Value = '  -1.08%  '
$ws.Range("D7").Value = '0.4581'
$ws.Range("E7").Value = '  -0.74%  '
$ws.Range("D8").Value = '0.3675'
$ws.Range("E8").Value = '  -1.05%  '
$ws.Range("D9").Value = '0.07153'
$ws.Range("E9").Value = '  -2.36%  '
$ws.Range("D10").Value = '0.8717'
$ws.Range("E10").Value = '  -1.19%  '
$ws.Range("D11").Value = '0.07766'
$ws.Range("E11").Value = '  -0.48%  '
$ws.Range("D12").Value = '19.54'
$ws.Range("E12").Value = '  -2.04%  '
$ws.Range("D13").Value = '1.812.69'
$ws.Range("E13").Value = '  -4.67%  '
$ws.Range("D14").Value = '5.313'
$ws.Range("E14").Value = '  -1.31%  '
$ws.Range("D15").Value = '6.377'
$ws.Range("E15").Value = '  -2.62%  '
$ws.Range("D16").Value = '86.87'
$ws.Range("E16").Value = '  -5.45%  '
$ws.Range("D17").Value = '1.007'
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("D18").Value = '0.000008703'
$ws.Range("E18").Value = '  -3.86%  '
$ws.Range("E19").Value = '  +0.52%  '
$ws.Range("D20").Value = '26.914.01'
$ws.Range("E20").Value = '  -1.69%  '
$ws.Range("D21").Value = '14.44'
$ws.Range("E21").Value = '  -2.48%  '
$ws.Range("D22").Value = '4.980'
$ws.Range("E22").Value = '  -2.92%  '
$ws.Range("D23").Value = '10.46'
$ws.Range("E23").Value = '  -0.65%  '
$ws.Range("D24").Value = '2.007'
$ws.Range("E24").Value = '  +4.13%  '
$ws.Range("D25").Value = '151.13'
$ws.Range("E25").Value = '  -0.68%  '
$ws.Range("D26").Value = '18.19'
$ws.Range("E26").Value = '  -0.92%  '
$ws.Range("D27").Value = '1.962'
$ws.Range("E27").Value = '  -5.48%  '
$ws.Range("D28").Value = '113.51'
$ws.Range("E28").Value = '  -2.15%  '
$ws.Range("D29").Value = '4.920'
$ws.Range("E29").Value = '  -3.64%  '
$ws.Range("D30").Value = '0.08796'
$ws.Range("E30").Value = '  -0.74%  '
$ws.Range("D31").Value = '3.080'
$ws.Range("E31").Value = '  +1.29%  '
$ws.Range("D32").Value = '0.7451'
$ws.Range("E32").Value = '  -3.36%  '
$ws.Range("D33").Value = '4.473'
$ws.Range("E33").Value = '  -0.57%  '
$ws.Range("E34").Value = '  -4.22%  '
$ws.Range("D35").Value = '2.512'
$ws.Range("E35").Value = '  -5.70%  '
$ws.Range("E36").Value = '  +0.92%  '
$ws.Range("D37").Value = '0.01936'
$ws.Range("E37").Value = '  -1.17%  '
$ws.Range("D38").Value = '0.05113'
$ws.Range("E38").Value = '  -2.28%  '
$ws.Range("D39").Value = '2.899'
$ws.Range("E39").Value = '  -1.90%  '
$ws.Range("D40").Value = '6.923'
$ws.Range("E40").Value = '  -1.39%  '
$ws.Range("D41").Value = '0.4960'
$ws.Range("E41").Value = '  -3.50%  '
$ws.Range("D42").Value = '0.1597'
$ws.Range("E42").Value = '  -2.42%  '
$ws.Range("D43").Value = '8.273'
$ws.Range("E43").Value = '  -1.98%  '
$ws.Range("D44").Value = '0.4673'
$ws.Range("E44").Value = '  -3.06%  '
$ws.Range("D45").Value = '1.006'
$ws.Range("E45").Value = '  +0.56%  '
$ws.Range("D46").Value = '10.11'
$ws.Range("E46").Value = '  -2.00%  '
$ws.Range("D47").Value = '101.50'
$ws.Range("E47").Value = '  -1.48%  '
$ws.Range("D48").Value = '1.605'
$ws.Range("E48").Value = '  -2.86%  '
$ws.Range("D49").Value = '0.06080'
$ws.Range("E49").Value = '  -2.27%  '
$ws.Range("D50").Value = '64.45'
$ws.Range("E50").Value = '  -2.01%  '
$ws.Range("D51").Value = '36.69'
$ws.Range("E51").Value = '  -0.34%  '
